$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New User Story cells in column B (rows 4-6)
$ws.Range("B4").Value = "Como um jogador novo deste tipo de jogos eu quero um tutorial para poder perceber como se começa a jogar"
$ws.Range("B5").Value = "Como um jogador já com alguma experiência e horas no jogo gostava que houvessem casamentos e estes formacem alianças para expandir a paz e comércio"
$ws.Range("B6").Value = "Como um jogador já com alguma experiência gostava de ter um sistema de preços conforme a demanda para que não possa abusar de loopholes"

# Shift the "To do" column content down/around
$ws.Range("D4").Value = "Analisar o código dado"

# Shift the "Doing" column content
$ws.Range("E4").Value = "Dar ideias no servidor de discord e discuti-las"
$ws.Range("E5").Value = "Jogar o jogo para conhecer melhor o projeto"
$ws.Range("E6").Value = "Pesquisar como se joga o jogo"

# Remove the now-obsolete tail of the "Done" column (rows 9-12)
$ws.Range("G9:G12").ClearContents() | Out-Null

# Update selection to match the new active cell
$ws.Range("E6").Select() | Out-Null

# Narrow column G to its new width
$ws.Columns("G").ColumnWidth = 51
